$wb = $excel.ActiveWorkbook

# --- BoM sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("BoM")

# "Data:" field — revision date bumped from 12/11/2024 to 30/11/2024
# (leading apostrophe forces text entry, same as the original cell, so the
# cell keeps its "text-formatted date" style instead of Excel reformatting
# it as a number/date)
$ws.Range("C9").Value = "'30/11/2024"

# Updated supplier stock (I), unit price (J) and subtotal (K) figures for
# the BoM line items (rows 13-43). Only the cells that actually changed are
# written; everything else (designators, quantities, formulas, styles) is
# left untouched.
$ws.Range("I13").Value = 310523
$ws.Range("J13").Value = 0.015
$ws.Range("K13").Value = 0.15

$ws.Range("I14").Value = 752073
$ws.Range("J14").Value = 0.2
$ws.Range("K14").Value = 0.4

$ws.Range("I15").Value = 1069929
$ws.Range("J15").Value = 0.26
$ws.Range("K15").Value = 0.26

$ws.Range("I16").Value = 2065446

$ws.Range("I17").Value = 4011
$ws.Range("J17").Value = 0.51
$ws.Range("K17").Value = 0.51

$ws.Range("I18").Value = 591372

$ws.Range("I19").Value = 571402

$ws.Range("I20").Value = 70179

$ws.Range("I21").Value = 34078

$ws.Range("I22").Value = 373124

$ws.Range("I23").Value = 12486

$ws.Range("I24").Value = 3849

$ws.Range("I25").Value = 112312

$ws.Range("I26").Value = 196590

$ws.Range("I27").Value = 129802

$ws.Range("I28").Value = 237873

$ws.Range("I29").Value = 806449

$ws.Range("I30").Value = 1439925

$ws.Range("I31").Value = 1286489

$ws.Range("I32").Value = 117776

$ws.Range("I33").Value = 181075

$ws.Range("I34").Value = 29900

$ws.Range("I35").Value = 147953

$ws.Range("I37").Value = 30248

$ws.Range("I38").Value = 11141
$ws.Range("J38").Value = 0.15
$ws.Range("K38").Value = 0.15

$ws.Range("I39").Value = 18943
$ws.Range("J39").Value = 0.33
$ws.Range("K39").Value = 0.33

$ws.Range("I40").Value = 29348

$ws.Range("I41").Value = 92195

$ws.Range("I42").Value = 15236

$ws.Range("I43").Value = 39272

# Totals (BoM!K44) and every downstream formula on the "Impostos" sheet
# (B12/C12, B14/C14, B15/C15, B17/C17, B24/C24, B25/C25, B26/C26, ...)
# recompute automatically from the figures above.
